# Insert a new data row at row 444 (pushing existing rows 444:550 down to 445:551)
# and populate the new row with the latest weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 444; existing rows 444-550 shift down to 445-551.
$ws.Rows.Item(444).Insert()

# Populate the newly inserted row 444 with the new record's data.
$ws.Cells.Item(444, 1).Value = 6
$ws.Cells.Item(444, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(444, 3).Value = "Metropolitana"
$ws.Cells.Item(444, 4).Value = 44855
$ws.Cells.Item(444, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(444, 5).Value = 13
$ws.Cells.Item(444, 6).Value = 100112039
$ws.Cells.Item(444, 7).Value = "Ciboulette"
$ws.Cells.Item(444, 8).Value = "Sin especificar"
$ws.Cells.Item(444, 9).Value = "Primera"
$ws.Cells.Item(444, 10).Value = 1300
$ws.Cells.Item(444, 11).Value = 700
$ws.Cells.Item(444, 12).Value = 800
$ws.Cells.Item(444, 13).Value = 749
$ws.Cells.Item(444, 14).Value = "`$/docena de atados"
$ws.Cells.Item(444, 15).Value = "Región Metropolitana"
$ws.Cells.Item(444, 16).Value = 250
$ws.Cells.Item(444, 17).Value = 3
$ws.Cells.Item(444, 18).Value = "Hortaliza"
